$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 handoff/handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 04:50:13"
$wsZhCn.Range("H2").Value = "2016-03-12 04:50:30"

# de-de sheet: row 2 handoff/handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 04:50:15"
$wsDeDe.Range("H2").Value = "2016-03-12 04:50:37"
